# url factory retain url allready generated
#
# This script replays a new "comparison" run of the URL checker against
# the workbook. A few existing "date of check" / "result" cells are
# refreshed with the results of the latest pass, and a brand-new
# comparison column is appended to the "atart" sheet (and a
# previously-empty comparison column is filled in on the "BIs" sheet).
# Finally the active sheet/selection is moved from "BIs" to "atart".

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param(
        $Range,
        [string]$Text
    )
    # Force the cell to stay a plain text value (Excel would otherwise
    # silently turn strings such as "100%" or "21/01/2016 09:48" into a
    # numeric percentage / date). Clearing the formatting afterwards
    # drops the temporary text number-format again so the cell keeps
    # using the sheet's default style, just like the rest of the
    # worksheet.
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.ClearFormats()
}

# ---------------------------------------------------------------------
# Sheet "batnet"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("batnet")
Set-TextValue $ws1.Range("D2") "21/01/2016 09:48"
Set-TextValue $ws1.Range("D4") "100%"

# ---------------------------------------------------------------------
# Sheet "jeka"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("jeka")
Set-TextValue $ws2.Range("D2") "21/01/2016 09:49"

# ---------------------------------------------------------------------
# Sheet "divers"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("divers")
Set-TextValue $ws3.Range("D2") "21/01/2016 09:49"
Set-TextValue $ws3.Range("D8") "58%"

# ---------------------------------------------------------------------
# Sheet "atart" - gains a brand new "comp3" column (E)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("atart")
Set-TextValue $ws4.Range("E1") "comp3"
Set-TextValue $ws4.Range("E2") "21/01/2016 09:51"
Set-TextValue $ws4.Range("E4") "100%"
Set-TextValue $ws4.Range("E5") "100%"
Set-TextValue $ws4.Range("E6") "100%"
Set-TextValue $ws4.Range("E7") "100%"
Set-TextValue $ws4.Range("E8") "100%"
Set-TextValue $ws4.Range("E9") "100%"
Set-TextValue $ws4.Range("E10") "100%"
Set-TextValue $ws4.Range("E11") "55%"

# ---------------------------------------------------------------------
# Sheet "BIs" - fills in the previously empty "comp2" column (D)
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("BIs")
Set-TextValue $ws5.Range("D2") "21/01/2016 09:51"
Set-TextValue $ws5.Range("D4") "94%"
Set-TextValue $ws5.Range("D5") "96%"
Set-TextValue $ws5.Range("D6") "97%"
Set-TextValue $ws5.Range("D7") "100%"
Set-TextValue $ws5.Range("D8") "97%"
Set-TextValue $ws5.Range("D9") "98%"
Set-TextValue $ws5.Range("D10") "100%"
Set-TextValue $ws5.Range("D11") "99%"
Set-TextValue $ws5.Range("D12") "99%"
Set-TextValue $ws5.Range("D13") "100%"
Set-TextValue $ws5.Range("D14") "100%"
Set-TextValue $ws5.Range("D15") "98%"
Set-TextValue $ws5.Range("D16") "94%"
Set-TextValue $ws5.Range("D17") "99%"
Set-TextValue $ws5.Range("D18") "100%"
Set-TextValue $ws5.Range("D19") "100%"
Set-TextValue $ws5.Range("D20") "97%"
Set-TextValue $ws5.Range("D21") "100%"
Set-TextValue $ws5.Range("D22") "99%"
Set-TextValue $ws5.Range("D23") "100%"
Set-TextValue $ws5.Range("D24") "99%"
Set-TextValue $ws5.Range("D25") "100%"
Set-TextValue $ws5.Range("D26") "83%"
Set-TextValue $ws5.Range("D27") "100%"
Set-TextValue $ws5.Range("D28") "82%"
Set-TextValue $ws5.Range("D29") "95%"
Set-TextValue $ws5.Range("D30") "87%"
Set-TextValue $ws5.Range("D31") "72%"
Set-TextValue $ws5.Range("D32") "99%"

# ---------------------------------------------------------------------
# Move the active tab / selection from "BIs" to "atart"
# ---------------------------------------------------------------------
$ws4.Activate()
$ws4.Range("E8").Select()
